$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 333 ("nafo informado") is removed; every row below shifts up by one.
$ws.Rows.Item(333).Delete()
